# Add the "ODI Batting Extra" sheet: per-match batting extras that
# complement the "ODI Batting" sheet (one row per match, keyed by
# MATCH_CODE, same row order/count as "ODI Batting").

$wb = $excel.ActiveWorkbook

# Add the new sheet directly after the last existing tab ("ODI Batting"),
# so the final tab order is: Player Info, ODI Batting, ODI Batting Extra.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row, styled to match the other sheets' headers (bold, thin box
# border, centered horizontally, top-aligned vertically).
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$headerRange = $newSheet.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# All columns except BATTING_POSITION (B) hold text values (match codes,
# counts and percentages are stored as text, not numbers), so force text
# formatting before writing them to avoid Excel's automatic number/percent
# conversion. BATTING_POSITION stays a real number.
$newSheet.Range("A2:A18").NumberFormat = "@"
$newSheet.Range("C2:F18").NumberFormat = "@"

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$data = @(
    @("3705", 3,    "4", "0", "24.19%", "NO"),
    @("3707", "",   "",  "",  "",       "NO"),
    @("3709", "",   "",  "",  "",       "NO"),
    @("3711", 3,    "1", "0", "6.14%",  "NO"),
    @("3721", "",   "",  "",  "",       "NO"),
    @("3722", 3,    "0", "0", "0.56%",  "NO"),
    @("3725", 3,    "0", "0", "",       "NO"),
    @("3730", "",   "",  "",  "",       "NO"),
    @("3754", 2,    "0", "0", "0.62%",  "NO"),
    @("3759", 7,    "1", "0", "8.19%",  "NO"),
    @("3764", 6,    "0", "0", "",       "NO"),
    @("3773", 8,    "0", "0", "7.04%",  "NO"),
    @("3778", "",   "",  "",  "",       "NO"),
    @("3785", 3,    "0", "0", "5.41%",  "NO"),
    @("4040", 7,    "0", "0", "6.13%",  "NO"),
    @("4043", 7,    "0", "0", "6.67%",  "NO"),
    @("4046", "",   "",  "",  "",       "NO")
)

$row = 2
foreach ($record in $data) {
    for ($col = 1; $col -le $record.Length; $col++) {
        $newSheet.Cells.Item($row, $col).Value = $record[$col - 1]
    }
    $row++
}

# Restore the originally active sheet/tab (adding a sheet shouldn't change
# which tab the workbook opens on).
$wb.Worksheets.Item("Player Info").Activate()

